$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, $val)
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue $ws "D2" "57.686.88"
$ws.Range("E2").Value = "  +0.02%  "

Set-TextValue $ws "D3" "3.061.39"
$ws.Range("E3").Value = "  +1.65%  "

$ws.Range("E4").Value = "  +0.05%  "

Set-TextValue $ws "D5" "514.43"
$ws.Range("E5").Value = "  +0.85%  "

Set-TextValue $ws "D6" "140.74"
$ws.Range("E6").Value = "  +0.81%  "

Set-TextValue $ws "D7" "0.999"
$ws.Range("E7").Value = "  -0.12%  "

Set-TextValue $ws "D8" "0.435"
$ws.Range("E8").Value = "  -0.57%  "

Set-TextValue $ws "D9" "7.29"

Set-TextValue $ws "D10" "0.110"
$ws.Range("E10").Value = "  -0.50%  "

$ws.Range("E11").Value = "  +2.78%  "

Set-TextValue $ws "D12" "3.578.72"
$ws.Range("E12").Value = "  +1.46%  "

$ws.Range("E13").Value = "  -3.24%  "

Set-TextValue $ws "D14" "26.88"
$ws.Range("E14").Value = "  +1.44%  "

$ws.Range("E15").Value = "  +3.10%  "

Set-TextValue $ws "D16" "57.567.68"
$ws.Range("E16").Value = "  -0.11%  "

$ws.Range("E17").Value = "  +0.43%  "

Set-TextValue $ws "D18" "3.068.17"
$ws.Range("E18").Value = "  +1.89%  "

Set-TextValue $ws "D19" "13.37"
$ws.Range("E19").Value = "  +4.09%  "

$ws.Range("E20").Value = "  +2.97%  "

Set-TextValue $ws "D21" "330.36"
$ws.Range("E21").Value = "  -0.10%  "

$ws.Range("E22").Value = "  +0.59%  "

$ws.Range("E23").Value = "  +1.99%  "

Set-TextValue $ws "D24" "65.89"
$ws.Range("E24").Value = "  +1.97%  "

Set-TextValue $ws "D25" "3.179.59"
$ws.Range("E25").Value = "  +1.38%  "

$ws.Range("E26").Value = "  -2.30%  "

$ws.Range("E27").Value = "  -0.30%  "

Set-TextValue $ws "D28" "0.0₃0904"
$ws.Range("E28").Value = "  -2.06%  "

Set-TextValue $ws "D29" "6.75"
$ws.Range("E29").Value = "  -0.62%  "

Set-TextValue $ws "D30" "7.31"
$ws.Range("E30").Value = "  -0.28%  "

Set-TextValue $ws "D31" "1.80"
$ws.Range("E31").Value = "  -0.65%  "

Set-TextValue $ws "D32" "1.20"
$ws.Range("E32").Value = "  +1.12%  "

Set-TextValue $ws "D33" "20.88"
$ws.Range("E33").Value = "  +1.47%  "

Set-TextValue $ws "D34" "153.53"
$ws.Range("E34").Value = "  -0.88%  "

$ws.Range("E35").Value = "  -2.04%  "

$ws.Range("E36").Value = "  +0.60%  "

$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws "D37" "1.28"
$ws.Range("E37").Value = "  +0.13%  "

$ws.Range("B38").Value = "EnergySwap"
$ws.Range("C38").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws "D38" "25.55"
$ws.Range("E38").Value = "  +4.45%  "

Set-TextValue $ws "D39" "0.0681"
$ws.Range("E39").Value = "  +0.60%  "

Set-TextValue $ws "D40" "37.13"
$ws.Range("E40").Value = "  -1.19%  "

Set-TextValue $ws "D41" "3.88"
$ws.Range("E41").Value = "  +0.71%  "

$ws.Range("E42").Value = "  +2.91%  "

Set-TextValue $ws "D43" "1.00"
$ws.Range("E43").Value = "  +0.33%  "

$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws "D44" "2.208.05"
$ws.Range("E44").Value = "  -0.96%  "

$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws "D45" "1.39"
$ws.Range("E45").Value = "  -1.34%  "

$ws.Range("E46").Value = "  +1.78%  "

$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextValue $ws "D47" "0.961"
$ws.Range("E47").Value = "  -2.44%  "

$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws "D48" "0.0245"
$ws.Range("E48").Value = "  +2.26%  "

Set-TextValue $ws "D49" "20.13"
$ws.Range("E49").Value = "  +3.65%  "

$ws.Range("E50").Value = "  -4.94%  "

$ws.Range("B51").Value = "TheGraph"
$ws.Range("C51").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue $ws "D51" "0.185"
$ws.Range("E51").Value = "  +0.41%  "
